# literature-mapping.xlsx: map literature to taxonomy and create a matrix of
# local/global with specific and agnostic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$excel.UserName = "Johannes Allgaier"

# ---------------------------------------------------------------------------
# 1) Build the local/global x specific/agnostic matrix: headers already exist
#    in B1:G1 (Model Specific, Model Agnostic, Local, Global, Intrinsic,
#    Post-Hoc). Mark DeepLIFT (row 2) as Model Specific + Local.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "x"
$ws.Range("D2").Value = "x"

# ---------------------------------------------------------------------------
# 2) Give the whole table (A1:G17) a full grid border and center the data
#    cells (B2:G17).
# ---------------------------------------------------------------------------
$table = $ws.Range("A1:G17")
$table.Borders.LineStyle = 1
$table.Borders.Weight = 2

$data = $ws.Range("B2:G17")
$data.HorizontalAlignment = -4108
$data.VerticalAlignment = -4108

$ws.Range("A1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3) Conditional formatting: highlight "x" cells in green (Excel's built-in
#    "Green Fill with Dark Green Text" style).
# ---------------------------------------------------------------------------
$fc = $data.FormatConditions.Add(1, 3, '"x"')
$fc.Interior.Color = 13561798
$fc.Font.Color = 24832

# ---------------------------------------------------------------------------
# 4) Threaded comments explaining each taxonomy column.
# ---------------------------------------------------------------------------
$ws.Range("B1").AddCommentThreaded("Model-specific interpretation tools are limited to specific model classes. The interpretation of regression weights in a linear model is a model-specific interpretation, since – by definition – the interpretation of intrinsically interpretable models is always model-specific. Tools that only work for the interpretation of e.g. neural networks are model-specific")
$ws.Range("C1").AddCommentThreaded("Model-agnostic tools can be used on any machine learning model and are applied after the model has been trained (post hoc). These agnostic methods usually work by analyzing feature input and output pairs. By definition, these methods cannot have access to model internals such as weights or structural information")
$ws.Range("D1").AddCommentThreaded("the interpretation method explains an individual prediction")
$ws.Range("E1").AddCommentThreaded("The interpretation method explains the entire model behavior.")
$ws.Range("F1").AddCommentThreaded("Intrinsic interpretability can be achieved by designing self-`nexplanatory models which incorporate interpretability di-`nrectly into the model structures.")
$ws.Range("G1").AddCommentThreaded("Post-hoc global explanation aims to provide a global un-`nderstanding about what knowledge has been acquired by`nthese pre-trained models, and illuminate the parameters or`nlearned representations in an intuitive manner to humans.")

# ---------------------------------------------------------------------------
# 5) Cosmetic: move the active selection like the author left it.
# ---------------------------------------------------------------------------
$ws.Range("A9").Select() | Out-Null
